# Adds a new "2022-Q3" sheet (right after "总计") with fund-holding data,
# and updates the "总计" (totals) sheet so its date list/row for the new
# quarter is inserted at the top (other quarters shift down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

function Set-HeaderCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $text
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

function Set-IndexCell($ws, $row, $col, $num) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $num
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

function Set-TextNumCell($ws, $row, $col, $text) {
    # Column holds numeric-looking text that must stay as text (not be
    # auto-converted to a number), mirroring the source workbook.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
}

# Header row (row 1) -- code/name columns have no A1 header in the source.
Set-HeaderCell $q3 1 2 "基金代码"
Set-HeaderCell $q3 1 3 "基金名称"
Set-HeaderCell $q3 1 4 "基金规模"
Set-HeaderCell $q3 1 5 "股票总仓位"
Set-HeaderCell $q3 1 6 "仓位占比"
Set-HeaderCell $q3 1 7 "持有市值(亿元)"
Set-HeaderCell $q3 1 8 "仓位排名"

$q3Rows = @(
    @(0, "002207", "前海开源金银珠宝主题精选混合C", "6.72", "90.85", "7.98", "0.5363", 6),
    @(1, "001302", "前海开源金银珠宝主题精选混合A", "3.99", "90.85", "7.98", "0.3184", 6),
    @(2, "003304", "前海开源沪港深核心资源灵活配置混合A", "3.45", "90.59", "7.77", "0.2681", 7),
    @(3, "003305", "前海开源沪港深核心资源灵活配置混合C", "1.89", "90.59", "7.77", "0.1469", 7)
)

$r = 2
foreach ($row in $q3Rows) {
    Set-IndexCell   $q3 $r 1 $row[0]
    Set-TextNumCell $q3 $r 2 $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    Set-TextNumCell $q3 $r 4 $row[3]
    Set-TextNumCell $q3 $r 5 $row[4]
    Set-TextNumCell $q3 $r 6 $row[5]
    Set-TextNumCell $q3 $r 7 $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Rewrite the "总计" sheet: a new top data row for 2022-Q3, all other
#    quarters shift down one row, and 2020-Q4 (now at the bottom) keeps
#    its original values on the newly appended row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 4, 1.27),
    @(1, "2022-Q2", 5, 2.63),
    @(2, "2022-Q1", 4, 1.5),
    @(3, "2021-Q4", 4, 1.32),
    @(4, "2021-Q3", 4, 1.65),
    @(5, "2021-Q2", 4, 2.53),
    @(6, "2021-Q1", 2, 1.04),
    @(7, "2020-Q4", 2, 1.95)
)

$r = 2
foreach ($row in $summaryRows) {
    Set-IndexCell $summary $r 1 $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Keep the originally-selected tab ("2020-Q4", the last sheet) active,
# rather than leaving the newly-inserted sheet selected.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Output "done"
